$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C3: re-colour to match the "in-progress" (orange) fill used by the AI row ---
# Copy the fill/alignment formatting from B4 (style s="3") onto C3 (was style s="5")
# via a format-only paste so the existing style slot gets re-used rather than a
# new duplicate style being created.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- C19: mark "Scenario script (WX)" row as Done, same formatting as C13/C16/C17 ---
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C19").Value = "Done"

# --- New comment on C3 ---
$ws.Range("C3").AddComment("Tom:`nSpeler rood sein inrijder Rsd") | Out-Null

# --- Selection / view state ---
$ws.Range("F4").Select()
